$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "282.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.87%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.068"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.70%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06625"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.78%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.285"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.70%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.372"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.27%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.372"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.53%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9361"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.42%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1573"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.13%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06086"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "12.66%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07540"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.18%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02939"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.27%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09056"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.11%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001576"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.09%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04440"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.15%"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006335"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.19%"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006127"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.54%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.473"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.03%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.238"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.19%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3206"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.11%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1291"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.64%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.075"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.06%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1515"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.55%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001174"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.24%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004441"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.46%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001242"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.27%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001608"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-2.50%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04154"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.15%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006232"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.88%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1247"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-11.59%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002008"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01148"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.95%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005478"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.21%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "25.93%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01299"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.81%"
